$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset rows (2,4,5,6,7,8) got reordered chronologically by date;
# row 3 is unaffected. Apply the new values cell-by-cell per the diff.

$ws.Range("D2").Value = 44623
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 1800
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 1900
$ws.Range("N2").Value = "$/paquete"
$ws.Range("P2").Value = 1900
$ws.Range("Q2").Value = 1

$ws.Range("D4").Value = 45225
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = 1750
$ws.Range("N4").Value = "$/docena de matas"
$ws.Range("P4").Value = 292
$ws.Range("Q4").Value = 6

$ws.Range("D5").Value = 45218
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 180
$ws.Range("K5").Value = 1400
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 1444
$ws.Range("N5").Value = "$/docena de matas"
$ws.Range("P5").Value = 241
$ws.Range("Q5").Value = 6

$ws.Range("D6").Value = 44377
$ws.Range("I6").Value = "Segunda"
$ws.Range("J6").Value = 550
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2800
$ws.Range("M6").Value = 2364
$ws.Range("N6").Value = "$/docena de matas"
$ws.Range("P6").Value = 394
$ws.Range("Q6").Value = 6

$ws.Range("D7").Value = 44370
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 1200
$ws.Range("M7").Value = 1080
$ws.Range("N7").Value = "$/docena de matas"
$ws.Range("P7").Value = 180
$ws.Range("Q7").Value = 6

$ws.Range("D8").Value = 44267
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 120
$ws.Range("K8").Value = 1500
$ws.Range("L8").Value = 1800
$ws.Range("M8").Value = 1650
$ws.Range("N8").Value = "$/docena de matas"
$ws.Range("P8").Value = 275
$ws.Range("Q8").Value = 6
